$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header column in H1, matching the style/formatting of the
# existing header cells (e.g. G1) by copying the formatting across.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the corresponding data value in H2 (plain numeric cell, like F2/G2)
$ws.Range("H2").Value = 1
